{"js": "// Update a set of calculated values across the \"Geometry\", \"Check\",\n// \"Shear reinforcement strength\" and \"Shear strength\" tables\n// (progress on shear/flexure testing numbers).\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// Each entry: [tableIndex, rowIndex, colIndex, oldValue, newValue]\n// Indices are 0-based and were confirmed against table.values contents.\nconst edits = [\n  [1, 4, 2, \"8.04\", \"0.0\"],     // Geometry: Longitudinal tension rebar (As)\n  [3, 1, 4, \"42.0\", \"42.45\"],   // Check: Stirrup spacing along length (Max.)\n  [3, 2, 4, \"42.0\", \"42.45\"],   // Check: Stirrup spacing along width (Max.)\n  [4, 4, 2, \"56.0\", \"56.6\"],    // Shear reinforcement strength: Effective height (d)\n  [4, 6, 2, \"22.82\", \"22.58\"],  // Shear reinforcement strength: Required shear reinforcing (Asw,req)\n  [4, 8, 2, \"49.57\", \"50.1\"],   // Shear reinforcement strength: Shear rebar strength (VRd,s)\n  [5, 1, 2, \"0.0072\", \"0.0\"],   // Shear strength: Longitudinal reinforcement ratio (\u03c1l)\n  [5, 2, 2, \"1.6\", \"1.59\"],     // Shear strength: k value (k)\n  [5, 6, 2, \"453.6\", \"458.46\"], // Shear strength: Maximum shear strength (VRd,max)\n  [5, 7, 2, \"49.57\", \"50.1\"],   // Shear strength: Total shear strength (VRd)\n  [5, 9, 2, \"10.09\", \"9.98\"],   // Shear strength: Demand Capacity Ratio (DCR)\n];\n\n// Load current cell text first so we can sanity-check against the\n// expected previous value before overwriting it.\nconst cells = edits.map(([tableIndex, rowIndex, colIndex]) => {\n  const cell = tables.items[tableIndex].getCell(rowIndex, colIndex);\n  cell.load(\"value\");\n  return cell;\n});\nawait context.sync();\n\nedits.forEach(([tableIndex, rowIndex, colIndex, oldValue, newValue], i) => {\n  const cell = cells[i];\n  if (cell.value.trim() !== oldValue) {\n    throw new Error(\n      `Unexpected value in table ${tableIndex} cell (${rowIndex},${colIndex}): ` +\n      `expected \"${oldValue}\" but found \"${cell.value}\"`\n    );\n  }\n  cell.value = newValue;\n});\nawait context.sync();\n", "ps1": "# Update a set of calculated values across the \"Geometry\", \"Check\",\n# \"Shear reinforcement strength\" and \"Shear strength\" tables\n# (progress on shear/flexure testing numbers).\n\n$d = $word.ActiveDocument\n\n# Each entry: Table index, Row index, Column index (all 1-based, Word COM\n# style), expected old value, new value.\n$edits = @(\n    @{ Table = 2; Row = 5;  Col = 3; Old = \"8.04\";   New = \"0.0\" },    # Geometry: Longitudinal tension rebar (As)\n    @{ Table = 4; Row = 2;  Col = 5; Old = \"42.0\";   New = \"42.45\" },  # Check: Stirrup spacing along length (Max.)\n    @{ Table = 4; Row = 3;  Col = 5; Old = \"42.0\";   New = \"42.45\" },  # Check: Stirrup spacing along width (Max.)\n    @{ Table = 5; Row = 5;  Col = 3; Old = \"56.0\";   New = \"56.6\" },   # Shear reinforcement strength: Effective height (d)\n    @{ Table = 5; Row = 7;  Col = 3; Old = \"22.82\";  New = \"22.58\" },  # Shear reinforcement strength: Required shear reinforcing (Asw,req)\n    @{ Table = 5; Row = 9;  Col = 3; Old = \"49.57\";  New = \"50.1\" },   # Shear reinforcement strength: Shear rebar strength (VRd,s)\n    @{ Table = 6; Row = 2;  Col = 3; Old = \"0.0072\"; New = \"0.0\" },    # Shear strength: Longitudinal reinforcement ratio (\u03c1l)\n    @{ Table = 6; Row = 3;  Col = 3; Old = \"1.6\";    New = \"1.59\" },   # Shear strength: k value (k)\n    @{ Table = 6; Row = 7;  Col = 3; Old = \"453.6\";  New = \"458.46\" }, # Shear strength: Maximum shear strength (VRd,max)\n    @{ Table = 6; Row = 8;  Col = 3; Old = \"49.57\";  New = \"50.1\" },   # Shear strength: Total shear strength (VRd)\n    @{ Table = 6; Row = 10; Col = 3; Old = \"10.09\";  New = \"9.98\" }    # Shear strength: Demand Capacity Ratio (DCR)\n)\n\nforeach ($edit in $edits) {\n    $cell = $d.Tables.Item($edit.Table).Cell($edit.Row, $edit.Col)\n    # Cell text includes the trailing cell-mark (CR + BEL); trim it before comparing.\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $edit.Old) {\n        throw \"Unexpected value in table $($edit.Table) cell ($($edit.Row),$($edit.Col)): expected '$($edit.Old)' but found '$current'\"\n    }\n    $cell.Range.Text = $edit.New\n}\n"}
